# Districts New-Old.xlsx : "added enforcement by jurisdiction"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark D4 (row for "E10 but no writeup") as done instead of wip
$ws.Range("D4").Value = "done"

# Mark D30 (row for "E1 & E3, no E3 writeup") as done instead of wip
$ws.Range("D30").Value = "done"

# Row 31 (S2 / S2) - enforcement done
$ws.Range("D31").Value = "done"

# Row 32 (W1 / S5) - enforcement done
$ws.Range("D32").Value = "done"

# Row 33 (W2 / S8) - old district value merges S7 & S8, enforcement in progress
$ws.Range("C33").Value = "S7 & S8"
$ws.Range("D33").Value = "wip"

# Update the active selection to reflect the last-edited cell
$ws.Range("D33").Select()
